# "Tasks and DRI" workbook update: add a "Cost and Benefit Analysis" task
# as the new item #7, pushing "Compilation and tidying up of codes",
# "Readme" and "Presentation Package" down by one row (and renumbering
# the "No" column accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Capture the existing contents of rows 8-10 (B/C/D columns) before they
# get shifted down, using Value() so the getter is actually invoked.
$b8 = $ws.Cells.Item(8, 2).Value()
$c8 = $ws.Cells.Item(8, 3).Value()
$d8 = $ws.Cells.Item(8, 4).Value()

$b9 = $ws.Cells.Item(9, 2).Value()
$c9 = $ws.Cells.Item(9, 3).Value()
$d9 = $ws.Cells.Item(9, 4).Value()

$b10 = $ws.Cells.Item(10, 2).Value()
$c10 = $ws.Cells.Item(10, 3).Value()
$d10 = $ws.Cells.Item(10, 4).Value()

# Create row 11 with the same formatting as row 10 (matching border/fill),
# then fill it with what used to be in row 10.
$ws.Range("A10:D10").Copy()
$ws.Range("A11:D11").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = $b10
$ws.Cells.Item(11, 3).Value = $c10
$ws.Cells.Item(11, 4).Value = $d10

# Shift what used to be row 9 into row 10.
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = $b9
$ws.Cells.Item(10, 3).Value = $c9
$ws.Cells.Item(10, 4).Value = $d9

# Shift what used to be row 8 into row 9.
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = $b8
$ws.Cells.Item(9, 3).Value = $c8
$ws.Cells.Item(9, 4).Value = $d8

# Put the new "Cost and Benefit Analysis" task into row 8.
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "Cost and Benefit Analysis"
$ws.Cells.Item(8, 3).Value = "Joel"
$ws.Cells.Item(8, 4).Value = "Done"

# Match the workbook's saved selection state.
$ws.Range("D8").Select()
